# Weekly update: insert a new observation row at row 28 (pushing all
# subsequent "Perejil" price rows down by one), matching the new weekly
# data feed cadence described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 28..98 down to 29..99, leaving a blank row 28 to populate.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the latest observation.
$ws.Cells.Item(28, 1).Value = 8
$ws.Cells.Item(28, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28, 3).Value = "Coquimbo"
$ws.Cells.Item(28, 4).Value = 44525
$ws.Cells.Item(28, 5).Value = 4
$ws.Cells.Item(28, 6).Value = 100112044
$ws.Cells.Item(28, 7).Value = "Perejil"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 3200
$ws.Cells.Item(28, 11).Value = 1500
$ws.Cells.Item(28, 12).Value = 2000
$ws.Cells.Item(28, 13).Value = 1750
$ws.Cells.Item(28, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(28, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(28, 16).Value = 1167
$ws.Cells.Item(28, 17).Value = 1.5
$ws.Cells.Item(28, 18).Value = "Hortaliza"
